$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Step 1: bump the date on row 12 (column B) from 45470 to 45471
$ws.Range("B12").Value = 45471

# Step 2: insert a blank range A2:F2, shifting A2:F12 down to A3:F13
$ws.Range("A2:F2").Insert()

# Step 3: cut the row that is now at A13:F13 (originally row 12) into A2:F2
$ws.Range("A13:F13").Cut($ws.Range("A2:F2"))

# Step 4: remove the now-empty range left behind by the cut, shifting rows back up
$ws.Range("A13:F13").Delete()

# Step 5: fix up row height and the "All" literal (re-sync to the canonical string)
$ws.Range("A2:F2").RowHeight = 409.5
$ws.Range("E2").Value = $ws.Range("E3").Value()

Write-Output "done"
